$wb = $excel.ActiveWorkbook

function Set-RowStats {
    param($ws, $row, $vit, $emp, $der, $total, $sofrido, $saldo)
    $ws.Cells.Item($row, 4).Value = $vit
    $ws.Cells.Item($row, 5).Value = $emp
    $ws.Cells.Item($row, 6).Value = $der
    $ws.Cells.Item($row, 7).Value = $total
    $ws.Cells.Item($row, 8).Value = $sofrido
    $ws.Cells.Item($row, 9).Value = $saldo
}

# --- Sheet "Grupo I" ---
$ws = $wb.Worksheets.Item("Grupo I")
$ws.Cells.Item(2, 2).Value = "Dom Camillo68"
$ws.Cells.Item(4, 2).Value = "Analove10 ITAQUI GRANDE!!"
$ws.Cells.Item(5, 2).Value = "Super Vasco f.c"
$ws.Cells.Item(3, 3).Value = 12
$ws.Cells.Item(4, 3).Value = 12
Set-RowStats $ws 2 0 0 0 520.8798828125 0 0
Set-RowStats $ws 3 0 0 0 513.55029296875 0 0
Set-RowStats $ws 4 0 0 0 526.31005859375 0 0
Set-RowStats $ws 5 0 0 0 505.749755859375 0 0

# --- Sheet "Grupo J" ---
$ws = $wb.Worksheets.Item("Grupo J")
$ws.Cells.Item(3, 2).Value = "pura bucha /botafogo"
$ws.Cells.Item(4, 2).Value = "Texas Club 2025"
$ws.Cells.Item(5, 2).Value = "TEAM LOPES 99"
$ws.Cells.Item(3, 3).Value = 12
$ws.Cells.Item(4, 3).Value = 12
Set-RowStats $ws 2 0 0 0 522.12060546875 0 0
Set-RowStats $ws 3 0 0 0 466.9794921875 0 0
Set-RowStats $ws 4 0 0 0 501.6298828125 0 0
Set-RowStats $ws 5 0 0 0 518.08984375 0 0

# --- Sheet "Grupo K" ---
$ws = $wb.Worksheets.Item("Grupo K")
$ws.Cells.Item(3, 2).Value = "Lá do Itaqui"
$ws.Cells.Item(3, 3).Value = 12
$ws.Cells.Item(4, 3).Value = 18
Set-RowStats $ws 2 0 0 0 552.1396484375 0 0
Set-RowStats $ws 3 0 0 0 507.43017578125 0 0
Set-RowStats $ws 4 0 0 0 513.68994140625 0 0
Set-RowStats $ws 5 0 0 0 443.39990234375 0 0

# --- Sheet "Grupo L" ---
$ws = $wb.Worksheets.Item("Grupo L")
$ws.Cells.Item(2, 2).Value = "TORRESMO COM PINGA"
$ws.Cells.Item(3, 3).Value = 12
Set-RowStats $ws 2 0 0 0 514.140625 0 0
Set-RowStats $ws 3 0 0 0 500.619384765625 0 0
Set-RowStats $ws 4 0 0 0 504.9794921875 0 0
Set-RowStats $ws 5 0 0 0 469.429931640625 0 0
